# fix: old predictions cleanup, dropping odds reset, git push flow
# Update btts_yes (R) and btts_no (S) odds for a handful of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("R7").Value = 1.62
$ws.Range("S7").Value = 2.2

$ws.Range("R17").Value = 2
$ws.Range("S17").Value = 1.73

$ws.Range("R23").Value = 1.5
$ws.Range("S23").Value = 2.5

$ws.Range("R26").Value = 1.44
$ws.Range("S26").Value = 2.62
